$wb = $excel.ActiveWorkbook

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZhCn     = $wb.Worksheets.Item("zh-cn")
$wsDeDe     = $wb.Worksheets.Item("de-de")

# Latest HO Xliff Generate Date (row for 001d19e8...)
$wsOverview.Range("G3").Value = "2016-08-30 16:54:14"

# zh-cn sheet, row for 001d19e8... (Correspond Handoff / Handback Datetime)
$wsZhCn.Range("H3").Value = "2016-08-30 16:54:03"
$wsZhCn.Range("K3").Value = "2016-08-30 16:54:31"

# de-de sheet, row for 001d19e8... (Correspond Handoff / Handback Datetime)
$wsDeDe.Range("H3").Value = "2016-08-30 16:54:14"
$wsDeDe.Range("K3").Value = "2016-08-30 16:54:38"
